# Fix the Ebios RM matrix
# - bump library_version on the library_content sheet
# - simplify the probability/impact "name" labels (drop the "Vx -"/"Gx -" prefix)
#   on the spec sheet, fixing a couple of mistranslations along the way
# - the workbook was left with the "spec" tab active/selected

$wb = $excel.ActiveWorkbook

$wsLibrary = $wb.Worksheets.Item("library_content")
$wsSpec    = $wb.Worksheets.Item("spec")

# library_content!B2 : library_version 1 -> 2
$wsLibrary.Range("B2").Value = 2

# spec sheet: probability rows (2-5) and impact rows (6-9) - columns E (English
# name) and K (French name) lose their "V1/V2/V3/V4 -" / "G1/G2/G3/G4 -" prefix.
$wsSpec.Range("E2").Value = "Certain"
$wsSpec.Range("K2").Value = "Certain"

$wsSpec.Range("E3").Value = "Very likely"
$wsSpec.Range("K3").Value = "Très vraisemblable"

$wsSpec.Range("E4").Value = "Likely"
$wsSpec.Range("K4").Value = "Vraisemblable"

$wsSpec.Range("E5").Value = "Unlikely"
$wsSpec.Range("K5").Value = "Peu vraisemblable"

$wsSpec.Range("E6").Value = "Minor"
$wsSpec.Range("K6").Value = "Mineur"

$wsSpec.Range("E7").Value = "Significant"
$wsSpec.Range("K7").Value = "Significatif"

$wsSpec.Range("E8").Value = "Important"
$wsSpec.Range("K8").Value = "Important"

$wsSpec.Range("E9").Value = "Critical"
$wsSpec.Range("K9").Value = "Critique"

# Selections / active tab: the file was last left on the "spec" sheet, with
# "library_content" having had B3 selected and "spec" having had F12 selected.
$wsLibrary.Range("B3").Select() | Out-Null
$wsSpec.Activate() | Out-Null
$wsSpec.Range("F12").Select() | Out-Null
